$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1 & 2 & 3. "left-top" (x2) and "left-bottom" (x1) -> "left-middle"
# ------------------------------------------------------------------
$d.Content.Find.Execute("left-top", $true, $false, $false, $false, $false, `
    $true, 1, $false, "left-middle", 2) | Out-Null

$d.Content.Find.Execute("left-bottom", $true, $false, $false, $false, $false, `
    $true, 1, $false, "left-middle", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "ay destrempé en fort bon <m>" -> "ay destrempé en <m>fort bon "
#    (the "fort bon " text moves to after the <m> tag; the <m> run
#     also loses its explicit Bold=False property)
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("fort bon <m>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$target1 = $d.Range($rng.Start, $rng.End)

$xml1 = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="0000ff"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t>&lt;m&gt;</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">fort bon </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target1.InsertXML($xml1)

# ------------------------------------------------------------------
# 5. " en ay moulé une fort" -> " en ay moulé une <al>fort"
#    (split the trailing "fort" off into its own run, inserting a new
#     <al> tag run before it)
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(" en ay moulé une fort", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$fortStart = $rng2.End - 4
$fortEnd = $rng2.End
$fortRng = $d.Range($fortStart, $fortEnd)

$xml2 = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="0000ff"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t>&lt;al&gt;</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t>fort</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$fortRng.InsertXML($xml2)

# ------------------------------------------------------------------
# 6. "petite <al>le" -> "petite le"
#    (remove the now-redundant <al> tag run; the surrounding "petite "
#     and "le" runs share identical formatting and merge automatically)
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("petite <al>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$alStart = $rng3.End - 4
$alEnd = $rng3.End
$alRng = $d.Range($alStart, $alEnd)
$alRng.Delete()
